$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.432.63'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '3.757.30'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.35'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.54'
$ws.Range("E6").Value = '  +1.57%  '
$ws.Range("D7").Value = '3.752.88'
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("E10").Value = '  +2.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.49'
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("E13").Value = '  +7.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.64'
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").Value = '4.388.26'
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").Value = '3.759.85'
$ws.Range("E16").Value = '  -0.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.70'
$ws.Range("E17").Value = '  +2.53%  '
$ws.Range("D18").Value = '67.480.39'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("E19").Value = '  -1.44%  '
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("E21").Value = '  -3.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '469.74'
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.720'
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000147'
$ws.Range("E24").Value = '  -7.50%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.94'
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.23'
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.18'
$ws.Range("E27").Value = '  +1.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.43'
$ws.Range("E28").Value = '  +4.69%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.91'
$ws.Range("E30").Value = '  -1.07%  '
$ws.Range("D31").Value = '3.905.97'
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.67'
$ws.Range("E32").Value = '  +1.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.64'
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("E35").Value = '  -2.82%  '
$ws.Range("D36").Value = '3.721.94'
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("E37").Value = '  +5.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.104'
$ws.Range("E38").Value = '  +1.14%  '
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.44%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.87'
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '400.45'
$ws.Range("E48").Value = '  -3.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000271'
$ws.Range("E49").Value = '  -7.45%  '
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '39.63'
$ws.Range("E50").Value = '  +4.53%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '140.09'
$ws.Range("E51").Value = '  -1.29%  '
